# Commit: "case with 380 kV done"
# Update Code/Results/Cases/Case_1_254/res_bus/vm_pu.xlsx values for rows 2-25
# in columns B,C,D,E,F,I,J,K,L,M,N (G stays 1, H/A untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colNums = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

$data = @(
    @(2, 1.02, 1.04409022963235, 1.062233826183949, 1.053260877141, 1.067103302838565, 1.049281608218312, 1.049157179334852, 1.0649557273581, 1.056007335212382, 1.069812064420158, 1.050647102869518),
    @(3, 1.02, 1.045004751672497, 1.06295056693624, 1.054063582760123, 1.067964173043777, 1.049538932090137, 1.04971898792739, 1.065487491057094, 1.056623076651251, 1.070488542877538, 1.051209709294697),
    @(4, 1.02, 1.045597069051292, 1.063414569487869, 1.054583818122376, 1.068522005684063, 1.04970429772942, 1.050082429337597, 1.065831126795069, 1.057021677186666, 1.07092641525513, 1.051573666833394),
    @(5, 1.02, 1.045846212205914, 1.063609687948267, 1.054802722466563, 1.068756706435065, 1.049773543742843, 1.050235198345429, 1.065975481998093, 1.057189289417313, 1.071110530223201, 1.051726652790741),
    @(6, 1.02, 1.045888052196653, 1.063642452184172, 1.054839488998702, 1.068796124693245, 1.049785154381261, 1.050260847626118, 1.06599971342674, 1.057217434590925, 1.071141445854652, 1.051752338496351),
    @(7, 1.02, 1.045600397593554, 1.063417176468916, 1.054586742359285, 1.068525141030779, 1.049705224074444, 1.050084470731436, 1.065833056107321, 1.057023916671664, 1.070928875278823, 1.051575711126247),
    @(8, 1.02, 1.044399179837729, 1.062476004816964, 1.053531982217231, 1.067394072927274, 1.04936880763854, 1.049347062391831, 1.065135531973051, 1.056215391002584, 1.070040652477849, 1.050837255582218),
    @(9, 1.02, 1.042286828379538, 1.060819325136547, 1.051679795565532, 1.065407128129297, 1.04876730040947, 1.048047042049896, 1.063903009419831, 1.054792058138073, 1.068476668630111, 1.049535389062283),
    @(10, 1.02, 1.040881590785776, 1.059716172113234, 1.050449425151043, 1.064086729050283, 1.048360491976998, 1.047180014800866, 1.063079121946469, 1.053844179754404, 1.067434889765073, 1.04866713053511),
    @(11, 1.02, 1.040273833403788, 1.059238824024302, 1.049917729390659, 1.063516006433246, 1.048182972223007, 1.046804513656587, 1.062721862609995, 1.05343399300514, 1.066984013809382, 1.048291096136136),
    @(12, 1.02, 1.040048194580396, 1.059061566115118, 1.049720395156814, 1.063304169424756, 1.048116828525955, 1.046665025970103, 1.062589085155234, 1.053281670298888, 1.06681657304777, 1.048151410361138),
    @(13, 1.02, 1.040096589907434, 1.059099586229564, 1.049762716697799, 1.063349602150552, 1.048131025836583, 1.046694946980574, 1.062617569762749, 1.053314342255681, 1.066852488069104, 1.048181373862875),
    @(14, 1.02, 1.040255179791078, 1.059224170791869, 1.049901414382842, 1.063498492752658, 1.048177508941234, 1.046792983755505, 1.062710888715486, 1.053421401161739, 1.066970172388888, 1.048279549861271),
    @(15, 1.02, 1.04035290677535, 1.059300938184533, 1.049986892018002, 1.063590249727856, 1.0482061215791, 1.046853386175945, 1.062768375645384, 1.053487368880672, 1.067042686215335, 1.048340038060076),
    @(16, 1.02, 1.040921940907928, 1.059747859113025, 1.050484734573277, 1.064124627643301, 1.048372244596474, 1.047204934136217, 1.063102821423482, 1.053871407887989, 1.067464817697916, 1.048692085258775),
    @(17, 1.02, 1.041279074528018, 1.060028288855034, 1.050797303698427, 1.064460102783217, 1.04847608316133, 1.047425432194492, 1.063312474623072, 1.054112373437363, 1.06772966983662, 1.048912896449584),
    @(18, 1.02, 1.04148745389355, 1.060191889971198, 1.050979722414793, 1.064655877951057, 1.048536518300412, 1.047554037902474, 1.063434712433388, 1.054252948686041, 1.067884174852714, 1.049041684792423),
    @(19, 1.02, 1.041558517601882, 1.060247678963568, 1.051041939809463, 1.06472264882762, 1.048557102686014, 1.047597887892233, 1.063476383957408, 1.05430088526692, 1.067936860630807, 1.049085597054197),
    @(20, 1.02, 1.041240750286866, 1.059998198159187, 1.05076375738279, 1.064424099306569, 1.048464955928566, 1.047401775564085, 1.063289985914516, 1.054086517610603, 1.067701251501019, 1.048889206224049),
    @(21, 1.02, 1.040208475995495, 1.059187482337999, 1.049860566907704, 1.063454643894881, 1.048163826477458, 1.046764114631898, 1.062683410681741, 1.053389873898132, 1.066935516322214, 1.048250639740197),
    @(22, 1.02, 1.039560077312855, 1.058678045697631, 1.049293628905603, 1.062846005169843, 1.047973308984681, 1.046363135480389, 1.062301596978029, 1.052952092111059, 1.066454269361212, 1.0478490911523),
    @(23, 1.02, 1.039903745309674, 1.058948079406452, 1.049594084436404, 1.063168570539078, 1.048074418037858, 1.046575707105731, 1.062504044524345, 1.05318414670907, 1.066709367942495, 1.048061964653733),
    @(24, 1.02, 1.041258067147567, 1.060011794753218, 1.050778915202519, 1.064440367424704, 1.048469984254444, 1.047412464999059, 1.063300147746532, 1.054098200671369, 1.067714092458137, 1.048899910839247),
    @(25, 1.02, 1.042832399074353, 1.061247394704037, 1.052157857903202, 1.06592006270672, 1.048923830654778, 1.048383194387731, 1.064222040494087, 1.055159851154315, 1.068880847573111, 1.049872018774995)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $colNums.Length; $i++) {
        $c = $colNums[$i]
        $val = $entry[$i + 1]
        $ws.Cells.Item($r, $c).Value = $val
    }
}
